{"js": "const replacements = [\n  [\"2025-05-04 Sunday\", \"2025-05-05 Monday\"],\n  [\"572\u00d72=\", \"839\u00d75=\"],\n  [\"164\u00d74=\", \"407\u00d75=\"],\n  [\"989\u00d76=\", \"458\u00d77=\"],\n  [\"601\u00d77=\", \"379\u00d76=\"],\n  [\"851\u00d74=\", \"309\u00d78=\"],\n  [\"240\u00d79=\", \"291\u00d76=\"],\n  [\"458\u00d76=\", \"907\u00d73=\"],\n  [\"562\u00d72=\", \"276\u00d76=\"],\n  [\"172\u00d72=\", \"490\u00d77=\"],\n  [\"180\u00d75=\", \"366\u00d79=\"],\n  [\"160\u00d74=\", \"570\u00d76=\"],\n  [\"616\u00d75=\", \"573\u00d79=\"],\n  [\"247\u00d78=\", \"317\u00d73=\"],\n  [\"147\u00d78=\", \"951\u00d77=\"],\n  [\"110\u00d78=\", \"645\u00d73=\"],\n  [\"513\u00d75=\", \"947\u00d73=\"],\n  [\"593\u00d78=\", \"613\u00d78=\"],\n  [\"905\u00d76=\", \"847\u00d76=\"],\n  [\"948\u00d76=\", \"542\u00d75=\"],\n  [\"349\u00d78=\", \"152\u00d77=\"],\n  [\"680\u00d75=\", \"630\u00d78=\"],\n  [\"795\u00d74=\", \"201\u00d78=\"],\n  [\"212\u00d77=\", \"510\u00d74=\"],\n  [\"493\u00d78=\", \"976\u00d76=\"],\n  [\"775\u00d73=\", \"634\u00d78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-04 Sunday\", \"2025-05-05 Monday\"),\n    @(\"572\u00d72=\", \"839\u00d75=\"),\n    @(\"164\u00d74=\", \"407\u00d75=\"),\n    @(\"989\u00d76=\", \"458\u00d77=\"),\n    @(\"601\u00d77=\", \"379\u00d76=\"),\n    @(\"851\u00d74=\", \"309\u00d78=\"),\n    @(\"240\u00d79=\", \"291\u00d76=\"),\n    @(\"458\u00d76=\", \"907\u00d73=\"),\n    @(\"562\u00d72=\", \"276\u00d76=\"),\n    @(\"172\u00d72=\", \"490\u00d77=\"),\n    @(\"180\u00d75=\", \"366\u00d79=\"),\n    @(\"160\u00d74=\", \"570\u00d76=\"),\n    @(\"616\u00d75=\", \"573\u00d79=\"),\n    @(\"247\u00d78=\", \"317\u00d73=\"),\n    @(\"147\u00d78=\", \"951\u00d77=\"),\n    @(\"110\u00d78=\", \"645\u00d73=\"),\n    @(\"513\u00d75=\", \"947\u00d73=\"),\n    @(\"593\u00d78=\", \"613\u00d78=\"),\n    @(\"905\u00d76=\", \"847\u00d76=\"),\n    @(\"948\u00d76=\", \"542\u00d75=\"),\n    @(\"349\u00d78=\", \"152\u00d77=\"),\n    @(\"680\u00d75=\", \"630\u00d78=\"),\n    @(\"795\u00d74=\", \"201\u00d78=\"),\n    @(\"212\u00d77=\", \"510\u00d74=\"),\n    @(\"493\u00d78=\", \"976\u00d76=\"),\n    @(\"775\u00d73=\", \"634\u00d78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2)\n}\n"}
